# Restored from revision of admin on 06/18/2020 07:17:27 AM.TEST Author: admin. Type: SAVE.
# Sheet "Rules": cell C10 changes from 18 to 1 (numeric value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
